# [PV-350][WIP] Replace hard coding of visual height with calculated value
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-01")

# Update header row labels
$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "Task Name"
$ws.Range("E1").Value = "Start"
$ws.Range("F1").Value = "Finish"

# Move active cell selection from F1 to F2
$ws.Activate()
$ws.Range("F2").Select()
